# Refresh the cryptos price-ticker rows (crypto price/volume(1h) columns D/E,
# plus a B/C/D/E rewrite for rows 42-43 where Kaspa and ARBITRUM swapped rank).
#
# Several "Price" column values are plain digits-and-dots text (e.g. "0.520",
# "8.90", "1.64") that Excel would otherwise auto-convert to a number on
# assignment (dropping trailing zeros / introducing float noise). Prefixing the
# string with a leading apostrophe forces Excel to keep it as literal text, the
# same trick a human typing into the grid would use; the apostrophe itself is
# not stored as part of the cell text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.902.16'
$ws.Range("E2").Value = '  +0.73%  '

# Row 3
$ws.Range("D3").Value = '1.625.85'
$ws.Range("E3").Value = '  +1.17%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '''214.15'
$ws.Range("E5").Value = '  +0.78%  '

# Row 6
$ws.Range("D6").Value = '''0.520'
$ws.Range("E6").Value = '  +0.12%  '

# Row 7
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = '''29.55'
$ws.Range("E8").Value = '  +8.59%  '

# Row 9
$ws.Range("E9").Value = '  +2.53%  '

# Row 10
$ws.Range("D10").Value = '''0.0611'
$ws.Range("E10").Value = '  +1.69%  '

# Row 11
$ws.Range("D11").Value = '''0.0914'
$ws.Range("E11").Value = '  +0.59%  '

# Row 12
$ws.Range("D12").Value = '1.859.07'
$ws.Range("E12").Value = '  +1.14%  '

# Row 13
$ws.Range("D13").Value = '1.623.14'
$ws.Range("E13").Value = '  +0.81%  '

# Row 14
$ws.Range("D14").Value = '''0.573'
$ws.Range("E14").Value = '  +6.55%  '

# Row 15
$ws.Range("E15").Value = '  +4.58%  '

# Row 16
$ws.Range("D16").Value = '29.960.01'
$ws.Range("E16").Value = '  +1.03%  '

# Row 17
$ws.Range("D17").Value = '''8.90'
$ws.Range("E17").Value = '  +17.04%  '

# Row 18
$ws.Range("D18").Value = '''64.67'
$ws.Range("E18").Value = '  +1.89%  '

# Row 19
$ws.Range("D19").Value = '''242.73'
$ws.Range("E19").Value = '  +0.58%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").Value = '  +1.79%  '

# Row 21
$ws.Range("E21").Value = '  -0.08%  '

# Row 22
$ws.Range("E22").Value = '  +3.62%  '

# Row 23
$ws.Range("D23").Value = '''9.63'
$ws.Range("E23").Value = '  +4.21%  '

# Row 24
$ws.Range("E24").Value = '  +1.96%  '

# Row 25
$ws.Range("D25").Value = '''157.72'
$ws.Range("E25").Value = '  +1.57%  '

# Row 26
$ws.Range("D26").Value = '''15.67'

# Row 27
$ws.Range("E27").Value = '  +2.43%  '

# Row 28
$ws.Range("E28").Value = '  +2.93%  '

# Row 29
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").Value = '''0.0489'
$ws.Range("E30").Value = '  +2.98%  '

# Row 31
$ws.Range("E31").Value = '  +4.85%  '

# Row 32
$ws.Range("E32").Value = '  +3.80%  '

# Row 33
$ws.Range("E33").Value = '  +3.10%  '

# Row 34
$ws.Range("D34").Value = '1.421.60'

# Row 35
$ws.Range("D35").Value = '''1.64'
$ws.Range("E35").Value = '  +6.37%  '

# Row 36
$ws.Range("E36").Value = '  +0.07%  '

# Row 37
$ws.Range("E37").Value = '  +1.49%  '

# Row 38
$ws.Range("E38").Value = '  -0.67%  '

# Row 39
$ws.Range("E39").Value = '  +3.28%  '

# Row 40
$ws.Range("D40").Value = '''0.557'
$ws.Range("E40").Value = '  +3.30%  '

# Row 41
$ws.Range("D41").Value = '''1.99'
$ws.Range("E41").Value = '  +0.99%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.0501'
$ws.Range("E42").Value = '  +2.58%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''0.832'
$ws.Range("E43").Value = '  +3.67%  '

# Row 44
$ws.Range("D44").Value = '''54.28'
$ws.Range("E44").Value = '  -0.19%  '

# Row 45
$ws.Range("D45").Value = '''69.37'
$ws.Range("E45").Value = '  +5.20%  '

# Row 46
$ws.Range("E46").Value = '  +7.23%  '

# Row 47
$ws.Range("D47").Value = '''0.998'
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("D48").Value = '''5.40'
$ws.Range("E48").Value = '  +1.96%  '

# Row 49
$ws.Range("D49").Value = '1.766.99'
$ws.Range("E49").Value = '  +1.07%  '

# Row 50
$ws.Range("D50").Value = '''88.83'
$ws.Range("E50").Value = '  +2.43%  '

# Row 51
$ws.Range("D51").Value = '0.0₆0108'
$ws.Range("E51").Value = '  +5.32%  '
